# Auto update Excel log
# Appends newly-logged sensor rows to the PIR, Humidity, Proximity and
# Camera sheets (mirrors the nightly SeniorConnect log-sync job).

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$r,
        [string]$date,
        [string]$time,
        [string]$hour,
        [string]$location,
        [string]$value,
        [string]$status
    )
    # Leading apostrophe forces each value to stay literal text instead of
    # being auto-converted (ISO dates -> date serials, "85.9%" -> 0.859, etc.)
    $ws.Cells.Item($r, 1).Value = "'" + $date
    $ws.Cells.Item($r, 2).Value = "'" + $time
    $ws.Cells.Item($r, 3).Value = "'" + $hour
    $ws.Cells.Item($r, 4).Value = "'" + $location
    $ws.Cells.Item($r, 5).Value = "'" + $value
    $ws.Cells.Item($r, 6).Value = "'" + $status
}

# ---------------------------------------------------------------------
# PIR sheet: append rows 194-207
# ---------------------------------------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")

$pirRows = @(
    @("2026-01-30", "18:37:23", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:23", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:27", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:32", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:37", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:42", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:47", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:52", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:37:57", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:38:02", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:38:07", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:38:12", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:38:17", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:38:22", "18:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 194
for ($i = 0; $i -lt $pirRows.Count; $i++) {
    $row = $pirRows[$i]
    $targetRow = $startRow + $i
    Add-LogRow $wsPIR $targetRow $row[0] $row[1] $row[2] $row[3] $row[4] $row[5]
}

# ---------------------------------------------------------------------
# Humidity sheet: append rows 127-136
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")

$humidityRows = @(
    @("2026-01-30", "18:37:23", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:37:27", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:37:32", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:37:38", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:37:48", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:37:53", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:37:57", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:38:07", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:38:12", "18:00", "Bathroom", "85.9%", "Active"),
    @("2026-01-30", "18:38:17", "18:00", "Bathroom", "85.8%", "Active")
)

$startRow = 127
for ($i = 0; $i -lt $humidityRows.Count; $i++) {
    $row = $humidityRows[$i]
    $targetRow = $startRow + $i
    Add-LogRow $wsHumidity $targetRow $row[0] $row[1] $row[2] $row[3] $row[4] $row[5]
}

# ---------------------------------------------------------------------
# Proximity sheet: append row 14 (single row - written directly so
# PowerShell doesn't unwrap a singleton array-of-arrays)
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $wsProximity 14 "2026-01-30" "18:37:37" "18:00" "Living Room Main Door" "EXIT" "User EXITED Living Room Main Door"

# ---------------------------------------------------------------------
# Camera sheet: append row 14 (single row)
# ---------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")
Add-LogRow $wsCamera 14 "2026-01-30" "18:37:37" "18:00" "Living Room Main Door" "Image Captured (EXIT)" "Active"

Write-Host "Appended log rows: PIR 194-207, Humidity 127-136, Proximity 14, Camera 14"
